# Auto-generated Excel COM-interop script
# Updates market-price-derived numeric cells (columns H-N) across several
# sheets to reflect the latest scheduled data pull. No formulas or
# structural changes are involved -- just literal value updates, plus two
# cells on CUL!132 (M132/N132) that are cleared because that row's HQ/NQ
# price data is no longer available (divide-by-zero-style N/A).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 280.23077
$ws.Range("I28").Value = 317.22223
$ws.Range("J28").Value = 197
$ws.Range("K28").Value = 317.22223
$ws.Range("L28").Value = 197
$ws.Range("M28").Value = 167.77777
$ws.Range("N28").Value = -1167

$ws.Range("H69").Value = 1646.3158
$ws.Range("J69").Value = 1610
$ws.Range("L69").Value = 4830
$ws.Range("N69").Value = -6578

$ws.Range("H72").Value = 1646.3158
$ws.Range("J72").Value = 1610
$ws.Range("L72").Value = 14490
$ws.Range("N72").Value = -23226

$ws.Range("H129").Value = 937.93335
$ws.Range("I129").Value = 474.625
$ws.Range("J129").Value = 993.2537
$ws.Range("K129").Value = 1423.875
$ws.Range("L129").Value = 2979.7611
$ws.Range("M129").Value = 3576.125
$ws.Range("N129").Value = -12979.7611

$ws.Range("H138").Value = 2947.0356
$ws.Range("J138").Value = 4116
$ws.Range("L138").Value = 12348
$ws.Range("N138").Value = -22628


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 19197.55
$ws.Range("I132").Value = 1394
$ws.Range("K132").Value = 4182
$ws.Range("M132").Value = -1652


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1238.4375
$ws.Range("I20").Value = 839.61536
$ws.Range("K20").Value = 839.61536
$ws.Range("M20").Value = -592.61536

$ws.Range("H86").Value = 1663.1923
$ws.Range("I86").Value = 1476.4117
$ws.Range("J86").Value = 2016
$ws.Range("K86").Value = 1476.4117
$ws.Range("L86").Value = 2016
$ws.Range("M86").Value = -353.4117000000001
$ws.Range("N86").Value = -4262

$ws.Range("H89").Value = 1663.1923
$ws.Range("I89").Value = 1476.4117
$ws.Range("J89").Value = 2016
$ws.Range("K89").Value = 7382.058500000001
$ws.Range("L89").Value = 10080
$ws.Range("M89").Value = -1766.058500000001
$ws.Range("N89").Value = -21312


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 678.76666
$ws.Range("I5").Value = 825.6667
$ws.Range("J5").Value = 642.0417
$ws.Range("K5").Value = 2477.0001
$ws.Range("L5").Value = 1926.1251
$ws.Range("M5").Value = -2365.0001
$ws.Range("N5").Value = -2150.1251

$ws.Range("H68").Value = 4606.613
$ws.Range("J68").Value = 5518.24
$ws.Range("L68").Value = 16554.72
$ws.Range("N68").Value = -18176.72

$ws.Range("H71").Value = 4606.613
$ws.Range("J71").Value = 5518.24
$ws.Range("L71").Value = 49664.16
$ws.Range("N71").Value = -57776.16

$ws.Range("H92").Value = 17862300
$ws.Range("J92").Value = 11333.333
$ws.Range("L92").Value = 33999.999
$ws.Range("N92").Value = -36495.999

$ws.Range("H123").Value = 3998
$ws.Range("I123").Value = 1200
$ws.Range("J123").Value = 4697.5
$ws.Range("K123").Value = 3600
$ws.Range("L123").Value = 14092.5
$ws.Range("M123").Value = -1150
$ws.Range("N123").Value = -18992.5

$ws.Range("H131").Value = 110713.03
$ws.Range("I131").Value = 810
$ws.Range("J131").Value = 118470.89
$ws.Range("K131").Value = 2430
$ws.Range("L131").Value = 355412.67
$ws.Range("M131").Value = 2610
$ws.Range("N131").Value = -365492.67

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0

$ws.Range("H133").Value = 3079.3635
$ws.Range("I133").Value = 2820
$ws.Range("J133").Value = 3227.5715
$ws.Range("K133").Value = 8460
$ws.Range("L133").Value = 9682.7145
$ws.Range("M133").Value = -3400
$ws.Range("N133").Value = -19802.7145

$ws.Range("H134").Value = 14298.625
$ws.Range("I134").Value = 14298.625
$ws.Range("K134").Value = 42895.875
$ws.Range("M134").Value = -37825.875

$ws.Range("H135").Value = 678.76666
$ws.Range("I135").Value = 825.6667
$ws.Range("J135").Value = 642.0417
$ws.Range("K135").Value = 7431.0003
$ws.Range("L135").Value = 5778.3753
$ws.Range("M135").Value = -4896.0003
$ws.Range("N135").Value = -10848.3753

$ws.Range("H136").Value = 2222.6316
$ws.Range("I136").Value = 949.2308
$ws.Range("J136").Value = 4981.6665
$ws.Range("K136").Value = 2847.6924
$ws.Range("L136").Value = 14944.9995
$ws.Range("M136").Value = 2252.3076
$ws.Range("N136").Value = -25144.9995

$ws.Range("H137").Value = 7641.8696
$ws.Range("I137").Value = 25757.25
$ws.Range("J137").Value = 3828.1052
$ws.Range("K137").Value = 77271.75
$ws.Range("L137").Value = 11484.3156
$ws.Range("M137").Value = -72171.75
$ws.Range("N137").Value = -21684.3156

$ws.Range("H138").Value = 2347.4
$ws.Range("I138").Value = 2226
$ws.Range("J138").Value = 2468.8
$ws.Range("K138").Value = 6678
$ws.Range("L138").Value = 7406.400000000001
$ws.Range("M138").Value = -1538
$ws.Range("N138").Value = -17686.4

$ws.Range("H139").Value = 28182.25
$ws.Range("I139").Value = 50514.5
$ws.Range("J139").Value = 5850
$ws.Range("K139").Value = 151543.5
$ws.Range("L139").Value = 17550
$ws.Range("M139").Value = -146403.5
$ws.Range("N139").Value = -27830

$ws.Range("H140").Value = 12374.5
$ws.Range("I140").Value = 16101.286
$ws.Range("J140").Value = 3678.6667
$ws.Range("K140").Value = 48303.858
$ws.Range("L140").Value = 11036.0001
$ws.Range("M140").Value = -43123.858
$ws.Range("N140").Value = -21396.0001

$ws.Range("H141").Value = 1798.3334
$ws.Range("I141").Value = 430
$ws.Range("K141").Value = 1290
$ws.Range("M141").Value = 3890


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1686.5217
$ws.Range("I102").Value = 1599.421
$ws.Range("K102").Value = 1599.421
$ws.Range("M102").Value = 22.57899999999995

$ws.Range("H132").Value = 61855.883
$ws.Range("I132").Value = 47663.22
$ws.Range("J132").Value = 170666.33
$ws.Range("K132").Value = 142989.66
$ws.Range("L132").Value = 511998.99
$ws.Range("M132").Value = -140459.66
$ws.Range("N132").Value = -517058.99


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 17473.656
$ws.Range("J136").Value = 3012.0908
$ws.Range("L136").Value = 9036.2724
$ws.Range("N136").Value = -14136.2724


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2692.6
$ws.Range("I81").Value = 2050
$ws.Range("J81").Value = 5263
$ws.Range("K81").Value = 4100
$ws.Range("L81").Value = 10526
$ws.Range("M81").Value = -3039
$ws.Range("N81").Value = -12648

$ws.Range("H84").Value = 2692.6
$ws.Range("I84").Value = 2050
$ws.Range("J84").Value = 5263
$ws.Range("K84").Value = 20500
$ws.Range("L84").Value = 52630
$ws.Range("M84").Value = -15196
$ws.Range("N84").Value = -63238

$ws.Range("H132").Value = 2885.1667
$ws.Range("I132").Value = 2425.25
$ws.Range("J132").Value = 3805
$ws.Range("K132").Value = 7275.75
$ws.Range("L132").Value = 11415
$ws.Range("M132").Value = -4745.75
$ws.Range("N132").Value = -16475


# ---- CUL row 132: HQ pricing data unavailable -> clear M/N entirely ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
